# Generate Report for handback
# Adds two new handback entries (12866827-... and fb54192e-...) as new rows
# to the Overview sheet and the two per-language detail sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet 1: Overview  (File Name | zh-cn | de-de)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewRows = @(
    @{ Row = 6; File = "12866827-d98c-4a17-b73d-265334bbbe8d.md" },
    @{ Row = 7; File = "fb54192e-400e-44ac-b23e-5224e823a2da.md" }
)

foreach ($r in $overviewRows) {
    $rowNum = $r.Row
    $fileName = $r.File

    $wsOverview.Cells.Item($rowNum, 1).Value = $fileName
    $wsOverview.Cells.Item($rowNum, 2).Value = $statusText
    $wsOverview.Cells.Item($rowNum, 3).Value = $statusText

    $wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($rowNum, 1), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$fileName", "", "", $fileName)
}

# ---------------------------------------------------------------------------
# Sheets 2 & 3: zh-cn / de-de detail sheets
# (Source File Name | Status | Correspond Handoff File | Correspond Handoff
#  Datetime | Target File | Correspond Handback File | Correspond Handback
#  DateTime | Handoff Reason | Dependency From)
# ---------------------------------------------------------------------------
$langSheets = @(
    @{
        Sheet = "zh-cn";
        Entries = @(
            @{
                Row = 6
                File = "12866827-d98c-4a17-b73d-265334bbbe8d.md"
                Xlf = "12866827-d98c-4a17-b73d-265334bbbe8d.79987c86519d0dbd026cdff0891d48d1f3e20ab6.zh-cn.xlf"
                HandoffDate = "2016-01-26 12:25:11"
                HandbackDate = "2016-01-26 12:26:01"
            },
            @{
                Row = 7
                File = "fb54192e-400e-44ac-b23e-5224e823a2da.md"
                Xlf = "fb54192e-400e-44ac-b23e-5224e823a2da.e2299bd3a4ac64c525b1f5ed1fd64c9f7c101ddf.zh-cn.xlf"
                HandoffDate = "2016-01-26 12:25:11"
                HandbackDate = "2016-01-26 12:26:01"
            }
        )
    },
    @{
        Sheet = "de-de";
        Entries = @(
            @{
                Row = 6
                File = "12866827-d98c-4a17-b73d-265334bbbe8d.md"
                Xlf = "12866827-d98c-4a17-b73d-265334bbbe8d.79987c86519d0dbd026cdff0891d48d1f3e20ab6.de-de.xlf"
                HandoffDate = "2016-01-26 12:25:22"
                HandbackDate = "2016-01-26 12:26:22"
            },
            @{
                Row = 7
                File = "fb54192e-400e-44ac-b23e-5224e823a2da.md"
                Xlf = "fb54192e-400e-44ac-b23e-5224e823a2da.e2299bd3a4ac64c525b1f5ed1fd64c9f7c101ddf.de-de.xlf"
                HandoffDate = "2016-01-26 12:25:22"
                HandbackDate = "2016-01-26 12:26:22"
            }
        )
    }
)

foreach ($langSheet in $langSheets) {
    $ws = $wb.Worksheets.Item($langSheet.Sheet)

    foreach ($e in $langSheet.Entries) {
        $rowNum = $e.Row
        $fileName = $e.File
        $xlfName = $e.Xlf

        $ws.Cells.Item($rowNum, 1).Value = $fileName
        $ws.Cells.Item($rowNum, 2).Value = $statusText
        $ws.Cells.Item($rowNum, 3).Value = $xlfName
        $ws.Cells.Item($rowNum, 4).Value = $e.HandoffDate
        $ws.Cells.Item($rowNum, 5).Value = $fileName
        $ws.Cells.Item($rowNum, 6).Value = $xlfName
        $ws.Cells.Item($rowNum, 7).Value = $e.HandbackDate
        $ws.Cells.Item($rowNum, 8).Value = "Include"

        $ws.Hyperlinks.Add($ws.Cells.Item($rowNum, 1), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$fileName", "", "", $fileName)
        $ws.Hyperlinks.Add($ws.Cells.Item($rowNum, 3), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.$($langSheet.Sheet)/xinjiang/ht/$xlfName", "", "", $xlfName)
        $ws.Hyperlinks.Add($ws.Cells.Item($rowNum, 5), "https://github.com/OpenLocalizationTestOrg/oltest.$($langSheet.Sheet)/blob/master/e2e/$fileName", "", "", $fileName)
        $ws.Hyperlinks.Add($ws.Cells.Item($rowNum, 6), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/master/ol-handback/OpenLocalizationTestOrg/oltest.$($langSheet.Sheet)/xinjiang/ht/$xlfName", "", "", $xlfName)
    }
}

Write-Output "handback rows added"
